# Apply the "Office Theme" colour scheme (previously used only by the
# Notes Master, ppt/theme/theme2.xml) to the presentation's main theme
# (ppt/theme/theme1.xml, used by the Slide Master), replacing the
# "Integral" / "Red Violet" colour scheme that lived there before.
#
# The 12 theme colours, in the fixed COM order
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), mapped to the
# RGB() long values PowerPoint's ThemeColorScheme exposes.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$newColors = @(
    0,          # dk1      -> 000000
    16777215,   # lt1      -> FFFFFF
    6968388,    # dk2      -> 44546A
    15132391,   # lt2      -> E7E6E6
    13998939,   # accent1  -> 5B9BD5
    3243501,    # accent2  -> ED7D31
    10855845,   # accent3  -> A5A5A5
    49407,      # accent4  -> FFC000
    12874308,   # accent5  -> 4472C4
    4697456,    # accent6  -> 70AD47
    12673797,   # hlink    -> 0563C1
    7491477     # folHlink -> 954F72
)

for ($i = 1; $i -le $newColors.Length; $i++) {
    $colorScheme.Item($i).RGB = $newColors[$i - 1]
}
